$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1142.25
$ws.Range("J29").Value = 2200
$ws.Range("L29").Value = 6600
$ws.Range("N29").Value = -7162
$ws.Range("H38").Value = 244.6875
$ws.Range("J38").Value = 916.6667
$ws.Range("L38").Value = 2750.0001
$ws.Range("N38").Value = -3494.0001
$ws.Range("H58").Value = 5682.8213
$ws.Range("J58").Value = 11911.923
$ws.Range("L58").Value = 35735.769
$ws.Range("N58").Value = -36035.769
$ws.Range("H111").Value = 1300
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1300
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 3900
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -10034
$ws.Range("H138").Value = 1812.3928
$ws.Range("I138").Value = 1414.0526
$ws.Range("J138").Value = 2653.3333
$ws.Range("K138").Value = 4242.1578
$ws.Range("L138").Value = 7959.999899999999
$ws.Range("M138").Value = 897.8422
$ws.Range("N138").Value = -18239.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14077.896
$ws.Range("I32").Value = 15604.803
$ws.Range("K32").Value = 15604.803
$ws.Range("M32").Value = -15317.803

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1080.5518
$ws.Range("I107").Value = 682.8333
$ws.Range("J107").Value = 1731.3636
$ws.Range("K107").Value = 682.8333
$ws.Range("L107").Value = 1731.3636
$ws.Range("M107").Value = 1237.1667
$ws.Range("N107").Value = -5571.3636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 47621076
$ws.Range("I94").Value = 333335550
$ws.Range("J94").Value = 1997.0555
$ws.Range("K94").Value = 333335550
$ws.Range("L94").Value = 1997.0555
$ws.Range("M94").Value = -333335099
$ws.Range("N94").Value = -2899.0555
$ws.Range("H107").Value = 545.0769
$ws.Range("I107").Value = 503.72223
$ws.Range("J107").Value = 638.125
$ws.Range("K107").Value = 503.72223
$ws.Range("L107").Value = 638.125
$ws.Range("M107").Value = 1416.27777
$ws.Range("N107").Value = -4478.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 927
$ws.Range("I44").Value = 867.1667
$ws.Range("J44").Value = 998.8
$ws.Range("K44").Value = 2601.5001
$ws.Range("L44").Value = 2996.4
$ws.Range("M44").Value = -2203.5001
$ws.Range("N44").Value = -3792.4
$ws.Range("H55").Value = 3340
$ws.Range("J55").Value = 3340
$ws.Range("L55").Value = 10020
$ws.Range("N55").Value = -10374
$ws.Range("H107").Value = 760.9048
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 825.26666
$ws.Range("K107").Value = 1800
$ws.Range("L107").Value = 2475.79998
$ws.Range("M107").Value = 120
$ws.Range("N107").Value = -6315.79998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4745.0586
$ws.Range("I80").Value = 3126.6
$ws.Range("J80").Value = 7057.143
$ws.Range("K80").Value = 3126.6
$ws.Range("L80").Value = 7057.143
$ws.Range("M80").Value = -2128.6
$ws.Range("N80").Value = -9053.143
$ws.Range("H83").Value = 4745.0586
$ws.Range("I83").Value = 3126.6
$ws.Range("J83").Value = 7057.143
$ws.Range("K83").Value = 15633
$ws.Range("L83").Value = 35285.715
$ws.Range("M83").Value = -10641
$ws.Range("N83").Value = -45269.715
$ws.Range("H122").Value = 47623628
$ws.Range("I122").Value = 76928630
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 230785890
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -230783440
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 58137
$ws.Range("I132").Value = 66642.19500000001
$ws.Range("J132").Value = 5404.8
$ws.Range("K132").Value = 199926.585
$ws.Range("L132").Value = 16214.4
$ws.Range("M132").Value = -197396.585
$ws.Range("N132").Value = -21274.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 540.1
$ws.Range("I22").Value = 433.33334
$ws.Range("J22").Value = 585.8570999999999
$ws.Range("K22").Value = 433.33334
$ws.Range("L22").Value = 585.8570999999999
$ws.Range("M22").Value = -138.33334
$ws.Range("N22").Value = -1175.8571
$ws.Range("H27").Value = 540.1
$ws.Range("I27").Value = 433.33334
$ws.Range("J27").Value = 585.8570999999999
$ws.Range("K27").Value = 433.33334
$ws.Range("L27").Value = 585.8570999999999
$ws.Range("M27").Value = -326.33334
$ws.Range("N27").Value = -799.8570999999999
$ws.Range("H61").Value = 33335696
$ws.Range("I61").Value = 2282.4
$ws.Range("J61").Value = 66669108
$ws.Range("K61").Value = 2282.4
$ws.Range("L61").Value = 66669108
$ws.Range("M61").Value = -2080.4
$ws.Range("N61").Value = -66669512
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H64").Value = 22400
$ws.Range("J64").Value = 22400
$ws.Range("L64").Value = 22400
$ws.Range("N64").Value = -22850
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H67").Value = 22400
$ws.Range("J67").Value = 22400
$ws.Range("L67").Value = 22400
$ws.Range("N67").Value = -23960
$ws.Range("H113").Value = 33335696
$ws.Range("I113").Value = 2282.4
$ws.Range("J113").Value = 66669108
$ws.Range("K113").Value = 2282.4
$ws.Range("L113").Value = 66669108
$ws.Range("M113").Value = -112.4000000000001
$ws.Range("N113").Value = -66673448
$ws.Range("H122").Value = 1925.0938
$ws.Range("I122").Value = 2154.4614
$ws.Range("J122").Value = 1768.1578
$ws.Range("K122").Value = 6463.3842
$ws.Range("L122").Value = 5304.4734
$ws.Range("M122").Value = -4013.3842
$ws.Range("N122").Value = -10204.4734

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 15000
$ws.Range("J63").Value = 15000
$ws.Range("L63").Value = 15000
$ws.Range("N63").Value = -16248
$ws.Range("H66").Value = 15000
$ws.Range("J66").Value = 15000
$ws.Range("L66").Value = 45000
$ws.Range("N66").Value = -51240
$ws.Range("H109").Value = 15000
$ws.Range("J109").Value = 15000
$ws.Range("L109").Value = 15000
$ws.Range("N109").Value = -17774
$ws.Range("H126").Value = 1522.4546
$ws.Range("I126").Value = 1766.2667
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 5298.800099999999
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -2828.800099999999
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 1014.6042
$ws.Range("I132").Value = 850.3095
$ws.Range("J132").Value = 2164.6667
$ws.Range("K132").Value = 2550.9285
$ws.Range("L132").Value = 6494.000100000001
$ws.Range("M132").Value = -20.92849999999999
$ws.Range("N132").Value = -11554.0001
